# Update gh-pages to output generated at 456a3b4
#
# This applies refreshed "want to go" counts (column F) scraped from
# bilibili show listings across the four sheets, and removes a duplicate
# row from the "全部类型" (all types) roll-up sheet, shifting the rows
# below it up by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: 展览 (Exhibitions)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1927
$ws1.Range("F3").Value  = 1537
$ws1.Range("F4").Value  = 905
$ws1.Range("F6").Value  = 13462
$ws1.Range("F7").Value  = 13304
$ws1.Range("F8").Value  = 1026
$ws1.Range("F13").Value = 6
$ws1.Range("F14").Value = 7
$ws1.Range("F15").Value = 699
$ws1.Range("F17").Value = 24
$ws1.Range("F22").Value = 410
$ws1.Range("F26").Value = 777

# ---------------------------------------------------------------
# Sheet: 演出 (Performances)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 100
$ws2.Range("F7").Value = 137
$ws2.Range("F8").Value = 588

# ---------------------------------------------------------------
# Sheet: 本地生活 (Local life)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 65

# ---------------------------------------------------------------
# Sheet: 全部类型 (All types roll-up)
# ---------------------------------------------------------------
# Row 29 was an accidental duplicate of row 28 (NIJISANJI EN). Delete it;
# Excel shifts rows 30:42 up into 29:41, shrinking the used range from
# A1:I42 to A1:I41.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Rows.Item(29).Delete()

# Refresh the "want to go" counts to match the other sheets / latest scrape.
$ws4.Range("F3").Value  = 1927
$ws4.Range("F4").Value  = 1537
$ws4.Range("F5").Value  = 905
$ws4.Range("F6").Value  = 100
$ws4.Range("F8").Value  = 13462
$ws4.Range("F9").Value  = 13304
$ws4.Range("F10").Value = 1026
$ws4.Range("F15").Value = 6
$ws4.Range("F16").Value = 7
$ws4.Range("F17").Value = 699
$ws4.Range("F21").Value = 24
$ws4.Range("F28").Value = 65
$ws4.Range("F29").Value = 410
$ws4.Range("F33").Value = 777
$ws4.Range("F34").Value = 137
$ws4.Range("F35").Value = 588
